$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.261.20'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '1.557.89'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3802'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3285'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.73'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.139'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07387'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.31'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.863'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.27%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.589.78'
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.770'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001078'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06666'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '86.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.445'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.32%  '
$ws.Range('D24').Value = '22.258.50'
$ws.Range('E24').Value = '  -1.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.299'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.573'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '150.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.32'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.946'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('D31').Value = '1.720.82'
$ws.Range('E31').Value = '  -2.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.083'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.933'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.920'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.444'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08225'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02363'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.19%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06361'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.370'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2165'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.243'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.06'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6089'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.755'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5914'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '123.44'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.978'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.179'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07087'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.77%  '
